# Applies the template-fix edit described in the commit message:
#   - Fix template functions - Add [v-megre] - Add [index:N]
#
# Concretely:
#  1. Widen the trailing paragraph's hanging indent from 324 -> 432 twips
#     (done FIRST, by paragraph index, before any table structural edit --
#     mutating Table.Rows height/indent invalidates previously-fetched
#     Paragraphs(N) handles in this COM host, so grab it up front).
#  2. Merge the ": " and "{{FileName}}" runs in the title paragraph into a
#     single run ": {{FileName}}" (keeping the en-US run formatting).
#  3. Widen the table indent (tblInd) from 432 -> 540 dxa.
#  4. Bump the header row height from 294 -> 299 dxa.
#  5. Grow the first data row height from 518 -> 728 dxa.
#  6. Append "[v-merge][index:{{Items.Column1}}]" markers to the
#     Items.Column1 / Items.Column2 template placeholders (Column1 must be
#     rewritten before Column2, otherwise the freshly-inserted
#     "{{Items.Column1}}" text inside Column2's replacement would itself
#     get matched and rewritten again).
#  7. Repoint the Items.Column3 / Items.Column4 placeholders at the new
#     Items.SubItems.Column1 / Items.SubItems.Column4 fields.

$d = $word.ActiveDocument

# 1. Trailing (last) paragraph hanging indent: 324 -> 432 twips, i.e.
#    16.2pt -> 21.6pt for both LeftIndent and (negative) FirstLineIndent.
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Format.LeftIndent = 21.6
$lastPara.Format.FirstLineIndent = -21.6

# 2. Merge ": " + "{{FileName}}" into a single run ": {{FileName}}".
$d.Content.Find.Execute(": {{FileName}}", $false, $false, $false, $false, `
    $false, $true, 1, $false, ": {{FileName}}", 2)

# 3. Table indent 432 -> 540 dxa (Word's Rows.LeftIndent is in points).
$table = $d.Tables.Item(1)
$table.Rows.LeftIndent = 27

# 4 & 5. Row heights (Row.Height is in points too: 299/20=14.95, 728/20=36.4).
$table.Rows.Item(1).Height = 14.95
$table.Rows.Item(2).Height = 36.4

# 6. Items.Column1, then Items.Column2, get the v-merge/index markers.
$d.Content.Find.Execute("{{Items.Column1}}", $false, $false, $false, $false, `
    $false, $true, 1, $false, `
    "{{Items.Column1}}[v-merge][index:{{Items.Column1}}]", 2)

$d.Content.Find.Execute("{{Items.Column2}}", $false, $false, $false, $false, `
    $false, $true, 1, $false, `
    "{{Items.Column2}}[v-merge][index:{{Items.Column1}}]", 2)

# 7. Column4 -> SubItems.Column4, Column3 -> SubItems.Column1.
$d.Content.Find.Execute("{{Items.Column4}}", $false, $false, $false, $false, `
    $false, $true, 1, $false, "{{Items.SubItems.Column4}}", 2)

$d.Content.Find.Execute("{{Items.Column3}}", $false, $false, $false, $false, `
    $false, $true, 1, $false, "{{Items.SubItems.Column1}}", 2)
